$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price update: insert two new rows of data at the top of this
# product's price block (rows 411-412), pushing the existing rows
# (old 411-434) down to 413-436.
$ws.Rows("411:412").Insert()

# --- Row 411 (new) ---
$ws.Cells.Item(411, 1).Value = 7
$ws.Cells.Item(411, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(411, 3).Value = "Ñuble"
$ws.Cells.Item(411, 4).Value = 44610
$ws.Cells.Item(411, 5).Value = 16
$ws.Cells.Item(411, 6).Value = "Fruta"
$ws.Cells.Item(411, 7).Value = 100106
$ws.Cells.Item(411, 8).Value = "Oleaginosos"
$ws.Cells.Item(411, 9).Value = 100106002
$ws.Cells.Item(411, 10).Value = "Palta"
$ws.Cells.Item(411, 11).Value = "Hass"
$ws.Cells.Item(411, 12).Value = "Primera"
$ws.Cells.Item(411, 13).Value = 600
$ws.Cells.Item(411, 14).Value = 2600
$ws.Cells.Item(411, 15).Value = 2800
$ws.Cells.Item(411, 16).Value = 2700
$ws.Cells.Item(411, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(411, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(411, 19).Value = 2700
$ws.Cells.Item(411, 20).Value = 1

# --- Row 412 (new) ---
$ws.Cells.Item(412, 1).Value = 7
$ws.Cells.Item(412, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(412, 3).Value = "Ñuble"
$ws.Cells.Item(412, 4).Value = 44610
$ws.Cells.Item(412, 5).Value = 16
$ws.Cells.Item(412, 6).Value = "Fruta"
$ws.Cells.Item(412, 7).Value = 100106
$ws.Cells.Item(412, 8).Value = "Oleaginosos"
$ws.Cells.Item(412, 9).Value = 100106002
$ws.Cells.Item(412, 10).Value = "Palta"
$ws.Cells.Item(412, 11).Value = "Hass"
$ws.Cells.Item(412, 12).Value = "Segunda"
$ws.Cells.Item(412, 13).Value = 500
$ws.Cells.Item(412, 14).Value = 2200
$ws.Cells.Item(412, 15).Value = 2400
$ws.Cells.Item(412, 16).Value = 2300
$ws.Cells.Item(412, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(412, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(412, 19).Value = 2300
$ws.Cells.Item(412, 20).Value = 1
